$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update card costs / creature stats (stone = G, wood = H columns)
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 3
$ws.Range("G5").Value = 2
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 1

# Update active selection to H5
$ws.Range("H5").Select()
